$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column E ("SampleName") rows 2-16 get a "Leak" suffix appended to record
# that these samples were run as part of the leakage test.
$ws.Range("E2").Value = "Mix1Leak"
$ws.Range("E3").Value = "Mix1Leak"
$ws.Range("E4").Value = "Mix1Leak"
$ws.Range("E5").Value = "Mix2Leak"
$ws.Range("E6").Value = "Mix2Leak"
$ws.Range("E7").Value = "Mix2Leak"
$ws.Range("E8").Value = "3N2OLeak"
$ws.Range("E9").Value = "3N2OLeak"
$ws.Range("E10").Value = "3N2OLeak"
$ws.Range("E11").Value = "10N2OLeak"
$ws.Range("E12").Value = "10N2OLeak"
$ws.Range("E13").Value = "10N2OLeak"
$ws.Range("E14").Value = "3KCO2Leak"
$ws.Range("E15").Value = "3KCO2Leak"
$ws.Range("E16").Value = "3KCO2Leak"

# Leave the cursor where the author last left it while entering this data.
$ws.Range("E16").Select() | Out-Null
